$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in row 2 with the dataset's sample values. The id-like values in
# columns C:H are numeric-looking but must stay text, so they're entered
# with a leading apostrophe (quote-prefix), same as typing them in Excel.
$ws.Range("A2").Value = "GenCor2019"
$ws.Range("B2").Value = "AR-X"
$ws.Range("C2").Value = "'32"
$ws.Range("D2").Value = "'1"
$ws.Range("E2").Value = "'4"
$ws.Range("F2").Value = "'1"
$ws.Range("G2").Value = "'9009"
$ws.Range("H2").Value = "'26"

# Drop the quote-prefix formatting flag picked up above so the cells keep
# the sheet's default (unstyled) look, matching row 2 in the target sheet.
$ws.Range("C2:H2").ClearFormats()

# Match the updated page margins (Excel's modern defaults in inches).
$ws.PageSetup.LeftMargin = 0.7 * 72
$ws.PageSetup.RightMargin = 0.7 * 72
$ws.PageSetup.TopMargin = 0.75 * 72
$ws.PageSetup.BottomMargin = 0.75 * 72
$ws.PageSetup.HeaderMargin = 0.3 * 72
$ws.PageSetup.FooterMargin = 0.3 * 72
